$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 52
$ws1.Range("F9").Value = 1982
$ws1.Range("F10").Value = 351
$ws1.Range("F11").Value = 4670
$ws1.Range("F12").Value = 82
$ws1.Range("F13").Value = 328

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 52
$ws4.Range("F13").Value = 1982
$ws4.Range("F14").Value = 351
$ws4.Range("F15").Value = 4670
$ws4.Range("F16").Value = 82
$ws4.Range("F17").Value = 328
